$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 10:51:15"
$wsZh.Range("H2").Value = "2016-03-19 10:51:33"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 10:51:18"
$wsDe.Range("H2").Value = "2016-03-19 10:51:38"
